# Auto-generated Excel COM-interop edit script
# Applies updated market-price data to the Moogle_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1391.7333
$ws.Range("I33").Value = 757.0417
$ws.Range("J33").Value = 3930.5
$ws.Range("K33").Value = 757.0417
$ws.Range("L33").Value = 3930.5
$ws.Range("M33").Value = -528.0417
$ws.Range("N33").Value = -4388.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6749.9707
$ws.Range("I43").Value = 7048.355
$ws.Range("K43").Value = 7048.355
$ws.Range("M43").Value = -6979.355

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7189.9
$ws.Range("I76").Value = 6999.6665
$ws.Range("J76").Value = 7271.4287
$ws.Range("K76").Value = 6999.6665
$ws.Range("L76").Value = 7271.4287
$ws.Range("M76").Value = -6684.6665
$ws.Range("N76").Value = -7901.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 7189.9
$ws.Range("I79").Value = 6999.6665
$ws.Range("J79").Value = 7271.4287
$ws.Range("K79").Value = 6999.6665
$ws.Range("L79").Value = 7271.4287
$ws.Range("M79").Value = -5907.6665
$ws.Range("N79").Value = -9455.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3649.5
$ws.Range("I80").Value = 2034
$ws.Range("J80").Value = 4090.0908
$ws.Range("K80").Value = 6102
$ws.Range("L80").Value = 12270.2724
$ws.Range("M80").Value = -5104
$ws.Range("N80").Value = -14266.2724

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 3649.5
$ws.Range("I83").Value = 2034
$ws.Range("J83").Value = 4090.0908
$ws.Range("K83").Value = 18306
$ws.Range("L83").Value = 36810.8172
$ws.Range("M83").Value = -13314
$ws.Range("N83").Value = -46794.8172

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8801.25
$ws.Range("I116").Value = 4343.4287
$ws.Range("J116").Value = 40006
$ws.Range("K116").Value = 4343.4287
$ws.Range("L116").Value = 40006
$ws.Range("M116").Value = -901.4287000000004
$ws.Range("N116").Value = -46890

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2742.394
$ws.Range("I132").Value = 2200.7083
$ws.Range("J132").Value = 4186.8887
$ws.Range("K132").Value = 6602.124899999999
$ws.Range("L132").Value = 12560.6661
$ws.Range("M132").Value = -4072.124899999999
$ws.Range("N132").Value = -17620.6661

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2102.0881
$ws.Range("I138").Value = 1799.2963
$ws.Range("J138").Value = 3270
$ws.Range("K138").Value = 5397.8889
$ws.Range("L138").Value = 9810
$ws.Range("M138").Value = -257.8888999999999
$ws.Range("N138").Value = -20090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 71249.75
$ws.Range("J24").Value = 71249.75
$ws.Range("L24").Value = 71249.75
$ws.Range("N24").Value = -71997.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9143.431
$ws.Range("I32").Value = 5375.121
$ws.Range("J32").Value = 28274.846
$ws.Range("K32").Value = 5375.121
$ws.Range("L32").Value = 28274.846
$ws.Range("M32").Value = -5088.121
$ws.Range("N32").Value = -28848.846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 86500
$ws.Range("J58").Value = 86500
$ws.Range("L58").Value = 86500
$ws.Range("N58").Value = -87360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4466.25
$ws.Range("I61").Value = 3737
$ws.Range("J61").Value = 5924.75
$ws.Range("K61").Value = 3737
$ws.Range("L61").Value = 5924.75
$ws.Range("M61").Value = -3525
$ws.Range("N61").Value = -6348.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3676.1035
$ws.Range("I74").Value = 2279.6667
$ws.Range("J74").Value = 10379
$ws.Range("K74").Value = 2279.6667
$ws.Range("L74").Value = 10379
$ws.Range("M74").Value = -1405.6667
$ws.Range("N74").Value = -12127

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3676.1035
$ws.Range("I77").Value = 2279.6667
$ws.Range("J77").Value = 10379
$ws.Range("K77").Value = 11398.3335
$ws.Range("L77").Value = 51895
$ws.Range("M77").Value = -7030.333500000001
$ws.Range("N77").Value = -60631

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 71249.75
$ws.Range("J100").Value = 71249.75
$ws.Range("L100").Value = 71249.75
$ws.Range("N100").Value = -73413.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1655.875
$ws.Range("I110").Value = 1655.875
$ws.Range("K110").Value = 1655.875
$ws.Range("M110").Value = 389.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4466.25
$ws.Range("I136").Value = 3737
$ws.Range("J136").Value = 5924.75
$ws.Range("K136").Value = 11211
$ws.Range("L136").Value = 17774.25
$ws.Range("M136").Value = -8661
$ws.Range("N136").Value = -22874.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 102427.71
$ws.Range("J139").Value = 102427.71
$ws.Range("L139").Value = 102427.71
$ws.Range("N139").Value = -112707.71

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2732.2666
$ws.Range("I134").Value = 2502.224
$ws.Range("K134").Value = 7506.672
$ws.Range("M134").Value = -4971.672

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6717.0684
$ws.Range("I31").Value = 3038.7273
$ws.Range("K31").Value = 3038.7273
$ws.Range("M31").Value = -2743.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6717.0684
$ws.Range("I34").Value = 3038.7273
$ws.Range("K34").Value = 3038.7273
$ws.Range("M34").Value = -2836.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 12214.5
$ws.Range("J96").Value = 12214.5
$ws.Range("L96").Value = 12214.5
$ws.Range("N96").Value = -17706.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2221.1538
$ws.Range("I46").Value = 1059.8
$ws.Range("J46").Value = 2947
$ws.Range("K46").Value = 3179.4
$ws.Range("L46").Value = 8841
$ws.Range("M46").Value = -3088.4
$ws.Range("N46").Value = -9023

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 14550
$ws.Range("I86").Value = 25100
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 75300
$ws.Range("L86").Value = 12000
$ws.Range("M86").Value = -74114
$ws.Range("N86").Value = -14372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 14550
$ws.Range("I89").Value = 25100
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 225900
$ws.Range("L89").Value = 36000
$ws.Range("M89").Value = -219972
$ws.Range("N89").Value = -47856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1942.6957
$ws.Range("I102").Value = 1176.5264
$ws.Range("J102").Value = 5582
$ws.Range("K102").Value = 1176.5264
$ws.Range("L102").Value = 5582
$ws.Range("M102").Value = 445.4736
$ws.Range("N102").Value = -8826

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4038.1428
$ws.Range("I132").Value = 3302.8857
$ws.Range("K132").Value = 9908.6571
$ws.Range("M132").Value = -7378.6571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5090.0713
$ws.Range("I7").Value = 5090.0713
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5090.0713
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4978.0713
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5207.9375
$ws.Range("I122").Value = 4059.0715
$ws.Range("J122").Value = 13250
$ws.Range("K122").Value = 12177.2145
$ws.Range("L122").Value = 39750
$ws.Range("M122").Value = -9727.2145
$ws.Range("N122").Value = -44650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5090.0713
$ws.Range("I126").Value = 5090.0713
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15270.2139
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12800.2139
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 54999.5
$ws.Range("J53").Value = 54999.5
$ws.Range("L53").Value = 54999.5
$ws.Range("N53").Value = -56213.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 50000
$ws.Range("J58").Value = 50000
$ws.Range("L58").Value = 50000
$ws.Range("N58").Value = -50616

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2933.6897
$ws.Range("I132").Value = 2387.5
$ws.Range("J132").Value = 4147.4443
$ws.Range("K132").Value = 7162.5
$ws.Range("L132").Value = 12442.3329
$ws.Range("M132").Value = -4632.5
$ws.Range("N132").Value = -17502.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1862.6061
$ws.Range("I136").Value = 1810.0435
$ws.Range("J136").Value = 1983.5
$ws.Range("K136").Value = 5430.1305
$ws.Range("L136").Value = 5950.5
$ws.Range("M136").Value = -2880.1305
$ws.Range("N136").Value = -11050.5
